$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert three new columns (B:D) on the UKA/MBA/LAU/MYD sheets, holding
#    the new low-end TEMP readings (14, 16, 18) that precede the existing
#    20..32 series in row 1. The existing B:H data shifts right to E:K.
#    Rows 2-4 are all the placeholder value 43000 regardless of column, so
#    after the shift the newly inserted B:D cells there simply need that
#    same value (and the row's normal border style) filled back in.
# ---------------------------------------------------------------------------
$sheetNames = @("UKA", "MBA", "LAU", "MYD")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Shift existing columns B..H right by 3 (insert 3 blank columns at B).
    $ws.Range("B1:D1").EntireColumn.Insert(-4161)

    # Rows 2-4: restore the ordinary data-cell border/fill (copied from the
    # first post-shift data column, E) onto the newly inserted B:D cells,
    # then fill in the placeholder values.
    $ws.Range("E2:E4").Copy() | Out-Null
    $ws.Range("B2:D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("B2:D4").Value2 = 43000

    # Row 1: new header values.
    $ws.Cells.Item(1, 2).Value = 14
    $ws.Cells.Item(1, 3).Value = 16
    $ws.Cells.Item(1, 4).Value = 18

    # Give the three new row-1 header cells the "open" left edge look (no
    # left border, thin border on the other three sides) instead of the
    # column-A style (medium left border) that Insert() propagated to them.
    foreach ($c in 2..4) {
        $cell = $ws.Cells.Item(1, $c)
        $cell.Borders.Item(7).LineStyle = 0
        $cell.Borders.Item(8).LineStyle = 1
        $cell.Borders.Item(8).Weight = 2
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(9).Weight = 2
        $cell.Borders.Item(10).LineStyle = 1
        $cell.Borders.Item(10).Weight = 2
    }

    # Selection now spans the wider A1:K4 table.
    $ws.Range("A1:K4").Select()
}

# ---------------------------------------------------------------------------
# 2) MYD was the active/selected tab; clear its special "I2" selection in
#    favor of the same A1:K4 block selection used by the other coasted
#    sheets (done above), then move the active tab back to WIL. This must
#    happen LAST, since selecting a range on another sheet re-activates it.
# ---------------------------------------------------------------------------
$wil = $wb.Worksheets.Item("WIL")
$wil.Activate()
$wil.Select()
